# edit.ps1 - apply the "use +2.0f instead of rand()" data regeneration edit
# (re-generated the two random-data blocks, added a third labeled data block,
#  added 5 new chart series, and repositioned/resized the chart.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: label cell (shared string 'fix #1') added in front of the existing 0..16 header row ---
$ws.Range("A12").Value2 = "fix #1"

# Row 13: fix#2 series 1 - regenerated data
$row13 = New-Object 'object[,]' 1,17
$row13[0,0] = 344.71169200000003
$row13[0,1] = 338.59280799999999
$row13[0,2] = 325.934617
$row13[0,3] = 344.317047
$row13[0,4] = 325.27729900000003
$row13[0,5] = 318.01305400000001
$row13[0,6] = 331.39467500000001
$row13[0,7] = 339.19288999999998
$row13[0,8] = 330.74798600000003
$row13[0,9] = 328.77432900000002
$row13[0,10] = 324.89948399999997
$row13[0,11] = 331.76573999999999
$row13[0,12] = 330.33281099999999
$row13[0,13] = 321.31094899999999
$row13[0,14] = 325.50493899999998
$row13[0,15] = 342.43057199999998
$row13[0,16] = 328.49892399999999
$ws.Range("B13:R13").Value2 = $row13

# Row 14: fix#2 series 2 - regenerated data
$row14 = New-Object 'object[,]' 1,17
$row14[0,0] = 337.51856400000003
$row14[0,1] = 365.60731900000002
$row14[0,2] = 409.06487800000002
$row14[0,3] = 682.43081700000005
$row14[0,4] = 629.65353300000004
$row14[0,5] = 661.28818899999999
$row14[0,6] = 632.72117600000001
$row14[0,7] = 702.82712100000003
$row14[0,8] = 645.59863099999995
$row14[0,9] = 664.45182699999998
$row14[0,10] = 625.25400999999999
$row14[0,11] = 664.55117900000005
$row14[0,12] = 641.72495700000002
$row14[0,13] = 651.25366299999996
$row14[0,14] = 629.84190999999998
$row14[0,15] = 690.40509599999996
$row14[0,16] = 643.46958800000004
$ws.Range("B14:R14").Value2 = $row14

# Row 15: fix#2 series 3 - regenerated data
$row15 = New-Object 'object[,]' 1,17
$row15[0,0] = 333.73659800000001
$row15[0,1] = 372.41520600000001
$row15[0,2] = 413.052458
$row15[0,3] = 593.78896699999996
$row15[0,4] = 545.88132599999994
$row15[0,5] = 534.73076300000002
$row15[0,6] = 536.05650000000003
$row15[0,7] = 622.66500599999995
$row15[0,8] = 533.27644999999995
$row15[0,9] = 571.23271999999997
$row15[0,10] = 539.80378099999996
$row15[0,11] = 617.426874
$row15[0,12] = 586.579072
$row15[0,13] = 907.19404899999995
$row15[0,14] = 949.64507000000003
$row15[0,15] = 1065.5301030000001
$row15[0,16] = 982.19766800000002
$ws.Range("B15:R15").Value2 = $row15

# --- Row 17: new label cell (shared string "fix #2") ---
$ws.Range("A17").Value2 = "fix #2"

# --- Rows 18-20: new labeled data block (3 constant-valued series) ---
$ws.Range("A18").Value2 = 1
$ws.Range("A19").Value2 = 2
$ws.Range("A20").Value2 = 4

# Row 18 - constant series
$row18 = New-Object 'object[,]' 1,17
for ($i = 0; $i -lt 17; $i++) { $row18[0,$i] = 338.98017800000002 }
$ws.Range("B18:R18").Value2 = $row18

# Row 19 - constant series
$row19 = New-Object 'object[,]' 1,17
for ($i = 0; $i -lt 17; $i++) { $row19[0,$i] = 700.90592100000003 }
$ws.Range("B19:R19").Value2 = $row19

# Row 20 - constant series
$row20 = New-Object 'object[,]' 1,17
for ($i = 0; $i -lt 17; $i++) { $row20[0,$i] = 1055.5760789999999 }
$ws.Range("B20:R20").Value2 = $row20

# --- Chart: add 5 new series (rows 16, 17 [still empty - placeholders], 18, 19, 20) ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$sc = $chart.SeriesCollection()

$newSer1 = $sc.NewSeries()
$newSer1.Values = "=Sheet1!`$B`$16:`$R`$16"

$newSer2 = $sc.NewSeries()
$newSer2.Values = "=Sheet1!`$B`$17:`$R`$17"

$newSer3 = $sc.NewSeries()
$newSer3.Values = "=Sheet1!`$B`$18:`$R`$18"

$newSer4 = $sc.NewSeries()
$newSer4.Values = "=Sheet1!`$B`$19:`$R`$19"

$newSer5 = $sc.NewSeries()
$newSer5.Values = "=Sheet1!`$B`$20:`$R`$20"

# --- Reposition / resize the chart (moved from A16:J38-ish to G10:S38-ish) ---
$co.Left = 417.1875
$co.Top = 160.5
$co.Width = 688.75
$co.Height = 455.5

# --- sheet view: selection moves to T22, and the sheet no longer needs to
#     scroll to keep row 3 at the top (topLeftCell reverts to the default) ---
$ws.Range("T22").Select()
